$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add "Coding Ninja" label in column A of row 28 (same row as the
# "to Move the last element to Front in a Linked List." entry)
$ws.Range("A28").Value = "Coding Ninja"

# Move the active selection down to the next empty row, matching the
# author's cursor position after making the edit
$ws.Activate()
$ws.Range("A29:B29").Select()
$excel.ActiveCell = $ws.Range("B29")
